$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "28.121.37"
Set-TextValue $ws "E2" "  -1.69%  "

Set-TextValue $ws "D3" "1.837.44"
Set-TextValue $ws "E3" "  -0.69%  "

Set-TextValue $ws "D4" "1.007"
Set-TextValue $ws "E4" "  +0.54%  "

Set-TextValue $ws "D5" "325.69"
Set-TextValue $ws "E5" "  -2.90%  "

Set-TextValue $ws "D6" "1.005"
Set-TextValue $ws "E6" "  +0.41%  "

Set-TextValue $ws "D7" "0.4637"
Set-TextValue $ws "E7" "  -0.35%  "

Set-TextValue $ws "D8" "0.3858"
Set-TextValue $ws "E8" "  -1.13%  "

Set-TextValue $ws "D9" "0.07842"
Set-TextValue $ws "E9" "  -0.84%  "

Set-TextValue $ws "D10" "0.9618"
Set-TextValue $ws "E10" "  -1.58%  "

Set-TextValue $ws "D11" "22.04"
Set-TextValue $ws "E11" "  -1.06%  "

Set-TextValue $ws "B12" "WrappedEther"
Set-TextValue $ws "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D12" "1.868.65"
Set-TextValue $ws "E12" "  -1.77%  "

Set-TextValue $ws "B13" "Polkadot"
Set-TextValue $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D13" "5.686"
Set-TextValue $ws "E13" "  -2.18%  "

Set-TextValue $ws "B14" "Chainlink"
Set-TextValue $ws "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D14" "6.861"
Set-TextValue $ws "E14" "  -1.46%  "

Set-TextValue $ws "D15" "0.06895"
Set-TextValue $ws "E15" "  -0.16%  "

Set-TextValue $ws "D16" "88.49"
Set-TextValue $ws "E16" "  +0.80%  "

Set-TextValue $ws "D17" "1.007"
Set-TextValue $ws "E17" "  +0.53%  "

Set-TextValue $ws "D18" "0.000009941"
Set-TextValue $ws "E18" "  -0.83%  "

Set-TextValue $ws "D19" "16.69"
Set-TextValue $ws "E19" "  -2.10%  "

Set-TextValue $ws "D20" "1.004"
Set-TextValue $ws "E20" "  +0.33%  "

Set-TextValue $ws "D21" "28.122.94"
Set-TextValue $ws "E21" "  -1.74%  "

Set-TextValue $ws "D22" "5.295"
Set-TextValue $ws "E22" "  -1.66%  "

Set-TextValue $ws "D23" "11.02"
Set-TextValue $ws "E23" "  -1.83%  "

Set-TextValue $ws "D24" "2.098"
Set-TextValue $ws "E24" "  -2.48%  "

Set-TextValue $ws "D25" "2.002.64"
Set-TextValue $ws "E25" "  -5.38%  "

Set-TextValue $ws "D26" "154.40"
Set-TextValue $ws "E26" "  +0.64%  "

Set-TextValue $ws "D27" "19.13"
Set-TextValue $ws "E27" "  -1.23%  "

Set-TextValue $ws "D28" "5.731"
Set-TextValue $ws "E28" "  -5.30%  "

Set-TextValue $ws "D29" "1.967"
Set-TextValue $ws "E29" "  -2.06%  "

Set-TextValue $ws "D30" "118.65"
Set-TextValue $ws "E30" "  +0.92%  "

Set-TextValue $ws "D31" "0.09259"
Set-TextValue $ws "E31" "  -0.95%  "

Set-TextValue $ws "D32" "0.9302"
Set-TextValue $ws "E32" "  -3.76%  "

Set-TextValue $ws "D33" "5.280"
Set-TextValue $ws "E33" "  -1.56%  "

Set-TextValue $ws "D34" "1.325"
Set-TextValue $ws "E34" "  -1.69%  "

Set-TextValue $ws "D35" "3.332"
Set-TextValue $ws "E35" "  -3.80%  "

Set-TextValue $ws "D36" "0.05807"
Set-TextValue $ws "E36" "  -4.85%  "

Set-TextValue $ws "D37" "0.02109"
Set-TextValue $ws "E37" "  -4.18%  "

Set-TextValue $ws "D38" "1.144"
Set-TextValue $ws "E38" "  -1.76%  "

Set-TextValue $ws "D39" "7.762"
Set-TextValue $ws "E39" "  +1.31%  "

Set-TextValue $ws "D40" "0.5581"
Set-TextValue $ws "E40" "  -2.11%  "

Set-TextValue $ws "D41" "9.846"
Set-TextValue $ws "E41" "  -2.64%  "

Set-TextValue $ws "D42" "0.1761"
Set-TextValue $ws "E42" "  -1.74%  "

Set-TextValue $ws "D43" "0.07235"
Set-TextValue $ws "E43" "  +2.04%  "

Set-TextValue $ws "D44" "11.68"
Set-TextValue $ws "E44" "  -0.97%  "

Set-TextValue $ws "D45" "0.5261"
Set-TextValue $ws "E45" "  -2.18%  "

Set-TextValue $ws "D46" "1.133"
Set-TextValue $ws "E46" "  -9.20%  "

Set-TextValue $ws "D47" "2.124"
Set-TextValue $ws "E47" "  -12.42%  "

Set-TextValue $ws "D48" "1.832"
Set-TextValue $ws "E48" "  -3.68%  "

Set-TextValue $ws "D49" "113.78"
Set-TextValue $ws "E49" "  +0.66%  "

Set-TextValue $ws "D50" "1.005"
Set-TextValue $ws "E50" "  +0.45%  "

Set-TextValue $ws "D51" "2.324"
Set-TextValue $ws "E51" "  -0.80%  "
